$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddContactInfo")

# Update the FullName value in row 2 from "Excel Project" to "Prakash Mijar"
$ws.Range("A2").Value = "Prakash Mijar"

# Move the active selection to A2 on this sheet
$ws.Activate()
$ws.Range("A2").Select()
